$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Prepare styles for new columns by copying the existing header style (s=1) ----
$ws.Range("D1").Copy()
$ws.Range("E1:L1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2:I2").PasteSpecial(-4122)
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("E3:I3").PasteSpecial(-4122)
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("E4:I4").PasteSpecial(-4122)
$ws.Range("L4").PasteSpecial(-4122)

# J (BirthPlace) and K (BirthDate) data cells (rows 2-4) use a distinct style
# derived from the same base font, but with an explicit (applied) default alignment.
# The header cells J1/K1 keep the plain copied style (s=1) from above.
$ws.Range("D2").Copy()
$ws.Range("J2:J4").PasteSpecial(-4122)
$ws.Range("J2:J4").WrapText = $false

$ws.Range("D2").Copy()
$ws.Range("K2:K4").PasteSpecial(-4122)
$ws.Range("K2:K4").WrapText = $false
$ws.Range("K2:K4").NumberFormat = "yyyy-mm-dd"

# ---- Header row ----
$ws.Range("E1").Value = "ClassRoom"
$ws.Range("F1").Value = "Gender"
$ws.Range("G1").Value = "NIS"
$ws.Range("H1").Value = "NISN"
$ws.Range("I1").Value = "Phone"
$ws.Range("J1").Value = "BirthPlace"
$ws.Range("K1").Value = "BirthDate"
$ws.Range("L1").Value = "Address"

# ---- Row 2 : Bill Cipher ----
$ws.Range("E2").Value = "9bb3b722-706d-4d35-bb92-c725b1075357"
$ws.Range("F2").Value = "pria"
$ws.Range("G2").Value = 202110469
$ws.Range("H2").Value = 46491333
$ws.Range("I2").Value = 85727721692
$ws.Range("J2").Value = "Bandung"
$ws.Range("K2").Value = 38353
$ws.Range("L2").Value = "Bandung Jawa Barat"

# ---- Row 3 : Amorhpous Shape ----
$ws.Range("E3").Value = "9bb3b73d-79aa-4d4f-8d6d-adbe487f48f2"
$ws.Range("F3").Value = "pria"
$ws.Range("G3").Value = 202110468
$ws.Range("H3").Value = 46491334
$ws.Range("I3").Value = 85727721693
$ws.Range("J3").Value = "Jakarta"
$ws.Range("K3").Value = 38354
$ws.Range("L3").Value = "Jakarta Jawa Barat"

# ---- Row 4 : Zanthar ----
$ws.Range("E4").Value = "9bb3b74e-61f9-435a-91ab-0bc58312a929"
$ws.Range("F4").Value = "wanita"
$ws.Range("G4").Value = 202110467
$ws.Range("H4").Value = 46491335
$ws.Range("I4").Value = 85727721694
$ws.Range("J4").Value = "Sukabumi"
$ws.Range("K4").Value = 38355
$ws.Range("L4").Value = "Sukabumi Jawa Barat"

# ---- Column widths: extend uniform width through column F ----
$ws.Range("F1").EntireColumn.ColumnWidth = 11.8
